# Refresh the cryptocurrency price/volume snapshot in the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values such as "30.478.12" or "53.40" that
# are formatted as plain text (not numbers) in the source data. Force the
# cells to Text format before writing so Excel does not reinterpret them as
# numeric/date values, then restore the default style so formatting stays
# untouched.
$priceCells = "D2,D3,D4,D5,D6,D7,D8,D9,D10,D11,D12,D13,D14,D15,D16,D17,D18,D19,D20,D22,D23,D24,D25,D26,D27,D28,D29,D30,D31,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D44,D46,D47,D48,D49,D50,D51"
foreach ($addr in $priceCells.Split(",")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "30.454.44"
$ws.Range("E2").Value = "  -1.25%  "
# Row 3
$ws.Range("D3").Value = "2.105.54"
$ws.Range("E3").Value = "  -0.55%  "
# Row 4
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.30%  "
# Row 5
$ws.Range("D5").Value = "334.02"
$ws.Range("E5").Value = "  +0.03%  "
# Row 6
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.16%  "
# Row 7
$ws.Range("D7").Value = "0.5250"
$ws.Range("E7").Value = "  -1.38%  "
# Row 8
$ws.Range("D8").Value = "0.4515"
$ws.Range("E8").Value = "  +2.03%  "
# Row 9
$ws.Range("D9").Value = "53.50"
$ws.Range("E9").Value = "  +12.78%  "
# Row 10
$ws.Range("D10").Value = "0.08989"
$ws.Range("E10").Value = "  -0.37%  "
# Row 11
$ws.Range("D11").Value = "1.179"
$ws.Range("E11").Value = "  -0.07%  "
# Row 12
$ws.Range("D12").Value = "24.40"
$ws.Range("E12").Value = "  -2.54%  "
# Row 13
$ws.Range("D13").Value = "2.102.54"
$ws.Range("E13").Value = "  -0.56%  "
# Row 14
$ws.Range("D14").Value = "6.784"
$ws.Range("E14").Value = "  +0.12%  "
# Row 15
$ws.Range("D15").Value = "7.804"
$ws.Range("E15").Value = "  -0.32%  "
# Row 16
$ws.Range("D16").Value = "96.56"
$ws.Range("E16").Value = "  -0.29%  "
# Row 17
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  +0.22%  "
# Row 18
$ws.Range("D18").Value = "0.00001131"
$ws.Range("E18").Value = "  -0.48%  "
# Row 19
$ws.Range("D19").Value = "0.06625"
$ws.Range("E19").Value = "  -0.93%  "
# Row 20
$ws.Range("D20").Value = "19.45"
$ws.Range("E20").Value = "  +1.45%  "
# Row 21
$ws.Range("E21").Value = "  +0.08%  "
# Row 22
$ws.Range("D22").Value = "6.320"
$ws.Range("E22").Value = "  -0.30%  "
# Row 23
$ws.Range("D23").Value = "30.514.40"
$ws.Range("E23").Value = "  -1.23%  "
# Row 24
$ws.Range("D24").Value = "12.36"
$ws.Range("E24").Value = "  +0.40%  "
# Row 25
$ws.Range("D25").Value = "2.348"
$ws.Range("E25").Value = "  +3.06%  "
# Row 26
$ws.Range("D26").Value = "2.349.72"
$ws.Range("E26").Value = "  -0.52%  "
# Row 27
$ws.Range("D27").Value = "22.34"
$ws.Range("E27").Value = "  -1.88%  "
# Row 28
$ws.Range("D28").Value = "2.577"
$ws.Range("E28").Value = "  -0.99%  "
# Row 29
$ws.Range("D29").Value = "163.55"
$ws.Range("E29").Value = "  +0.02%  "
# Row 30
$ws.Range("D30").Value = "132.93"
# Row 31
$ws.Range("D31").Value = "1.196"
$ws.Range("E31").Value = "  +0.15%  "
# Row 32
$ws.Range("D32").Value = "0.1074"
$ws.Range("E32").Value = "  -1.14%  "
# Row 33
$ws.Range("D33").Value = "1.658"
$ws.Range("E33").Value = "  +6.27%  "
# Row 34
$ws.Range("D34").Value = "6.155"
$ws.Range("E34").Value = "  -1.20%  "
# Row 35
$ws.Range("D35").Value = "3.927"
$ws.Range("E35").Value = "  -2.39%  "
# Row 36
$ws.Range("D36").Value = "10.54"
$ws.Range("E36").Value = "  +9.90%  "
# Row 37
$ws.Range("D37").Value = "0.02580"
$ws.Range("E37").Value = "  -1.09%  "
# Row 38
$ws.Range("D38").Value = "5.573"
$ws.Range("E38").Value = "  +0.24%  "
# Row 39
$ws.Range("D39").Value = "0.06818"
$ws.Range("E39").Value = "  +0.55%  "
# Row 40
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "12.75"
$ws.Range("E40").Value = "  -1.24%  "
# Row 41
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.2299"
$ws.Range("E41").Value = "  -0.45%  "
# Row 42
$ws.Range("D42").Value = "0.6897"
$ws.Range("E42").Value = "  +0.49%  "
# Row 43
$ws.Range("E43").Value = "  +0.29%  "
# Row 44
$ws.Range("D44").Value = "2.355"
$ws.Range("E44").Value = "  +5.08%  "
# Row 46
$ws.Range("D46").Value = "14.08"
$ws.Range("E46").Value = "  -0.25%  "
# Row 47
$ws.Range("D47").Value = "0.6382"
$ws.Range("E47").Value = "  -1.20%  "
# Row 48
$ws.Range("D48").Value = "3.657"
$ws.Range("E48").Value = "  -0.15%  "
# Row 49
$ws.Range("D49").Value = "1.246"
$ws.Range("E49").Value = "  -1.89%  "
# Row 50
$ws.Range("D50").Value = "1.228"
$ws.Range("E50").Value = "  +2.79%  "
# Row 51
$ws.Range("D51").Value = "83.48"
$ws.Range("E51").Value = "  +0.47%  "

foreach ($addr in $priceCells.Split(",")) {
    $ws.Range($addr).Style = "Normal"
}
